$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6220463333333334
$ws.Range("H2").Value = 1.866139
$ws.Range("I2").Value = 0.00505260120118785
$ws.Range("J2").Value = 0.00505260120118785
$ws.Range("M2").Value = 8.553891
$ws.Range("N2").Value = 25.661673
$ws.Range("O2").Value = 0.2062132866242743
$ws.Range("P2").Value = 0.2062132866242743
$ws.Range("Q2").Value = 5.320916532283
$ws.Range("R2").Value = 47.888248790547
$ws.Range("S2").Value = 0.001041913499698703
$ws.Range("T2").Value = 0.001041913499698703
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6220463333333334
$ws.Range("H3").Value = 1.866139
$ws.Range("I3").Value = 0.00505260120118785
$ws.Range("J3").Value = 0.00505260120118785
$ws.Range("O3").Value = 0.1165172631215423
$ws.Range("P3").Value = 0.1165172631215423
$ws.Range("Q3").Value = 3.006492170261556
$ws.Range("R3").Value = 27.058429532354
$ws.Range("S3").Value = 0.0005887152636070253
$ws.Range("T3").Value = 0.0005887152636070253
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6220463333333334
$ws.Range("H4").Value = 1.866139
$ws.Range("I4").Value = 0.00505260120118785
$ws.Range("J4").Value = 0.00505260120118785
$ws.Range("M4").Value = 15.27280066666667
$ws.Range("N4").Value = 45.818402
$ws.Range("O4").Value = 0.3681896836691911
$ws.Range("P4").Value = 0.3681896836691911
$ws.Range("Q4").Value = 9.50038965443089
$ws.Range("R4").Value = 85.503506889878
$ws.Range("S4").Value = 0.001860315637971929
$ws.Range("T4").Value = 0.001860315637971929
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6220463333333334
$ws.Range("H5").Value = 1.866139
$ws.Range("I5").Value = 0.00505260120118785
$ws.Range("J5").Value = 0.00505260120118785
$ws.Range("M5").Value = 12.82087433333333
$ws.Range("N5").Value = 38.46262299999999
$ws.Range("O5").Value = 0.3090797665849924
$ws.Range("P5").Value = 0.3090797665849924
$ws.Range("Q5").Value = 7.975177869177443
$ws.Range("R5").Value = 71.77660082259699
$ws.Range("S5").Value = 0.001561656799910193
$ws.Range("T5").Value = 0.001561656799910193
# Row 6
$ws.Range("G6").Value = 19.60581866666667
$ws.Range("H6").Value = 58.81745600000001
$ws.Range("I6").Value = 0.1592492032139157
$ws.Range("J6").Value = 0.1592492032139158
$ws.Range("M6").Value = 8.553891
$ws.Range("N6").Value = 25.661673
$ws.Range("O6").Value = 0.2062132866242743
$ws.Range("P6").Value = 0.2062132866242743
$ws.Range("Q6").Value = 167.706035840432
$ws.Range("R6").Value = 1509.354322563888
$ws.Range("S6").Value = 0.03283930158703851
$ws.Range("T6").Value = 0.03283930158703851
# Row 7
$ws.Range("G7").Value = 19.60581866666667
$ws.Range("H7").Value = 58.81745600000001
$ws.Range("I7").Value = 0.1592492032139157
$ws.Range("J7").Value = 0.1592492032139158
$ws.Range("O7").Value = 0.1165172631215423
$ws.Range("P7").Value = 0.1165172631215423
$ws.Range("Q7").Value = 94.75940481320178
$ws.Range("R7").Value = 852.8346433188161
$ws.Range("S7").Value = 0.01855528131277177
$ws.Range("T7").Value = 0.01855528131277178
# Row 8
$ws.Range("G8").Value = 19.60581866666667
$ws.Range("H8").Value = 58.81745600000001
$ws.Range("I8").Value = 0.1592492032139157
$ws.Range("J8").Value = 0.1592492032139158
$ws.Range("M8").Value = 15.27280066666667
$ws.Range("N8").Value = 45.818402
$ws.Range("O8").Value = 0.3681896836691911
$ws.Range("P8").Value = 0.3681896836691911
$ws.Range("Q8").Value = 299.4357604028125
$ws.Range("R8").Value = 2694.921843625312
$ws.Range("S8").Value = 0.05863391375590236
$ws.Range("T8").Value = 0.05863391375590237
# Row 9
$ws.Range("G9").Value = 19.60581866666667
$ws.Range("H9").Value = 58.81745600000001
$ws.Range("I9").Value = 0.1592492032139157
$ws.Range("J9").Value = 0.1592492032139158
$ws.Range("M9").Value = 12.82087433333333
$ws.Range("N9").Value = 38.46262299999999
$ws.Range("O9").Value = 0.3090797665849924
$ws.Range("P9").Value = 0.3090797665849924
$ws.Range("Q9").Value = 251.3637373274542
$ws.Range("R9").Value = 2262.273635947088
$ws.Range("S9").Value = 0.0492207065582031
$ws.Range("T9").Value = 0.0492207065582031
# Row 10
$ws.Range("G10").Value = 1.570446666666667
$ws.Range("H10").Value = 4.71134
$ws.Range("I10").Value = 0.01275602843261105
$ws.Range("J10").Value = 0.01275602843261106
$ws.Range("M10").Value = 8.553891
$ws.Range("N10").Value = 25.661673
$ws.Range("O10").Value = 0.2062132866242743
$ws.Range("P10").Value = 0.2062132866242743
$ws.Range("Q10").Value = 13.43342960798
$ws.Range("R10").Value = 120.90086647182
$ws.Range("S10").Value = 0.002630462547361415
$ws.Range("T10").Value = 0.002630462547361416
# Row 11
$ws.Range("G11").Value = 1.570446666666667
$ws.Range("H11").Value = 4.71134
$ws.Range("I11").Value = 0.01275602843261105
$ws.Range("J11").Value = 0.01275602843261106
$ws.Range("O11").Value = 0.1165172631215423
$ws.Range("P11").Value = 0.1165172631215423
$ws.Range("Q11").Value = 7.590327848804444
$ws.Range("R11").Value = 68.31295063924
$ws.Range("S11").Value = 0.001486297521268417
$ws.Range("T11").Value = 0.001486297521268417
# Row 12
$ws.Range("G12").Value = 1.570446666666667
$ws.Range("H12").Value = 4.71134
$ws.Range("I12").Value = 0.01275602843261105
$ws.Range("J12").Value = 0.01275602843261106
$ws.Range("M12").Value = 15.27280066666667
$ws.Range("N12").Value = 45.818402
$ws.Range("O12").Value = 0.3681896836691911
$ws.Range("P12").Value = 0.3681896836691911
$ws.Range("Q12").Value = 23.98511889763111
$ws.Range("R12").Value = 215.86607007868
$ws.Range("S12").Value = 0.004696638073478272
$ws.Range("T12").Value = 0.004696638073478272
# Row 13
$ws.Range("G13").Value = 1.570446666666667
$ws.Range("H13").Value = 4.71134
$ws.Range("I13").Value = 0.01275602843261105
$ws.Range("J13").Value = 0.01275602843261106
$ws.Range("M13").Value = 12.82087433333333
$ws.Range("N13").Value = 38.46262299999999
$ws.Range("O13").Value = 0.3090797665849924
$ws.Range("P13").Value = 0.3090797665849924
$ws.Range("Q13").Value = 20.13449936053555
$ws.Range("R13").Value = 181.21049424482
$ws.Range("S13").Value = 0.003942630290502951
$ws.Range("T13").Value = 0.003942630290502951
# Row 14
$ws.Range("G14").Value = 101.3157653333333
$ws.Range("H14").Value = 303.947296
$ws.Range("I14").Value = 0.8229421671522854
$ws.Range("J14").Value = 0.8229421671522854
$ws.Range("M14").Value = 8.553891
$ws.Range("N14").Value = 25.661673
$ws.Range("O14").Value = 0.2062132866242743
$ws.Range("P14").Value = 0.2062132866242743
$ws.Range("Q14").Value = 866.644013242912
$ws.Range("R14").Value = 7799.796119186208
$ws.Range("S14").Value = 0.1697016089901757
$ws.Range("T14").Value = 0.1697016089901757
# Row 15
$ws.Range("G15").Value = 101.3157653333333
$ws.Range("H15").Value = 303.947296
$ws.Range("I15").Value = 0.8229421671522854
$ws.Range("J15").Value = 0.8229421671522854
$ws.Range("O15").Value = 0.1165172631215423
$ws.Range("P15").Value = 0.1165172631215423
$ws.Range("Q15").Value = 489.6822613943395
$ws.Range("R15").Value = 4407.140352549056
$ws.Range("S15").Value = 0.09588696902389506
$ws.Range("T15").Value = 0.09588696902389506
# Row 16
$ws.Range("G16").Value = 101.3157653333333
$ws.Range("H16").Value = 303.947296
$ws.Range("I16").Value = 0.8229421671522854
$ws.Range("J16").Value = 0.8229421671522854
$ws.Range("M16").Value = 15.27280066666667
$ws.Range("N16").Value = 45.818402
$ws.Range("O16").Value = 0.3681896836691911
$ws.Range("P16").Value = 0.3681896836691911
$ws.Range("Q16").Value = 1547.375488326777
$ws.Range("R16").Value = 13926.37939494099
$ws.Range("S16").Value = 0.3029988162018386
$ws.Range("T16").Value = 0.3029988162018386
# Row 17
$ws.Range("G17").Value = 101.3157653333333
$ws.Range("H17").Value = 303.947296
$ws.Range("I17").Value = 0.8229421671522854
$ws.Range("J17").Value = 0.8229421671522854
$ws.Range("M17").Value = 12.82087433333333
$ws.Range("N17").Value = 38.46262299999999
$ws.Range("O17").Value = 0.3090797665849924
$ws.Range("P17").Value = 0.3090797665849924
$ws.Range("Q17").Value = 1298.956695324156
$ws.Range("R17").Value = 11690.61025791741
$ws.Range("S17").Value = 0.2543547729363762
$ws.Range("T17").Value = 0.2543547729363762

Write-Host "Applied all cell updates"
